$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 2 (room A1): capacity 20 -> 12, drop the stray "Loai phong" value in D2
$ws.Range("C2").Value = 12
$ws.Range("D2").ClearContents()

# --- New row 4: room "B1" (reuses the existing shared string), capacity 12
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "B1"
$ws.Range("C4").Value = 12

# --- New row 5: room "B2", capacity 12
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "B2"
$ws.Range("C5").Value = 12

# --- Row 3 now becomes room "A2" (was "B1") with capacity 12, D3 cleared
$ws.Range("B3").Value = "A2"
$ws.Range("C3").Value = 12
$ws.Range("D3").ClearContents()

# --- Update the active selection to reflect where the user left off editing
[void]$ws.Range("A6").Select()
